$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Estado de Cuenta" worker table lists CC + doc-number + name for each
# worker. The commit reorders/refreshes this small database: the worker
# previously shown in the first data row (YESIKA LUNA BASTIDAS /
# 1052041403) and the worker previously shown in the last data row
# (VIVIANA SAENZ LUNA / 1052052410) trade places, while the middle row
# (NATANAEL MENDOZA TORRES) stays put.

$ws.Range("C16").Value = "1052052410"
$ws.Range("D16").Value = "VIVIANA SAENZ LUNA"

$ws.Range("C18").Value = "1052041403"
$ws.Range("D18").Value = "YESIKA LUNA BASTIDAS"
